$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44162
$ws.Range("M2").Value = 300
# Row 3
$ws.Range("D3").Value = 44162
$ws.Range("M3").Value = 300
# Row 4
$ws.Range("D4").Value = 44165
$ws.Range("K4").Value = 'Castle Brite'
$ws.Range("N4").Value = 20500
$ws.Range("O4").Value = 21000
$ws.Range("P4").Value = 20750
$ws.Range("Q4").Value = '$/caja 15 kilos'
$ws.Range("S4").Value = 1383
$ws.Range("T4").Value = 15
# Row 5
$ws.Range("D5").Value = 44165
$ws.Range("K5").Value = 'Castle Brite'
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17750
$ws.Range("Q5").Value = '$/caja 15 kilos'
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1183
$ws.Range("T5").Value = 15
# Row 6
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 400
$ws.Range("N6").Value = 23500
$ws.Range("O6").Value = 24000
$ws.Range("P6").Value = 23750
$ws.Range("S6").Value = 1319
# Row 7
$ws.Range("D7").Value = 44189
$ws.Range("K7").Value = 'Dina'
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 200
$ws.Range("N7").Value = 21500
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 21750
$ws.Range("Q7").Value = '$/caja 18 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1208
$ws.Range("T7").Value = 18
# Row 8
$ws.Range("D8").Value = 44181
$ws.Range("K8").Value = 'Modesto'
$ws.Range("M8").Value = 16
$ws.Range("N8").Value = 495000
$ws.Range("O8").Value = 500000
$ws.Range("P8").Value = 497500
$ws.Range("Q8").Value = '$/bins (500 kilos)'
$ws.Range("S8").Value = 995
$ws.Range("T8").Value = 500
# Row 9
$ws.Range("D9").Value = 44181
$ws.Range("K9").Value = 'Modesto'
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 10
$ws.Range("N9").Value = 425000
$ws.Range("O9").Value = 430000
$ws.Range("P9").Value = 427500
$ws.Range("Q9").Value = '$/bins (500 kilos)'
$ws.Range("S9").Value = 855
$ws.Range("T9").Value = 500
# Row 10
$ws.Range("D10").Value = 44174
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 240
$ws.Range("N10").Value = 22500
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 22750
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("S10").Value = 1264
$ws.Range("T10").Value = 18
# Row 11
$ws.Range("D11").Value = 44187
$ws.Range("K11").Value = 'Dina'
$ws.Range("L11").Value = 'Especial'
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 22000
$ws.Range("O11").Value = 23000
$ws.Range("P11").Value = 22500
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("S11").Value = 1250
$ws.Range("T11").Value = 18
# Row 12
$ws.Range("D12").Value = 44186
$ws.Range("K12").Value = 'Dina'
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 22500
$ws.Range("O12").Value = 23000
$ws.Range("P12").Value = 22750
$ws.Range("Q12").Value = '$/caja 18 kilos'
$ws.Range("S12").Value = 1264
$ws.Range("T12").Value = 18
# Row 13
$ws.Range("D13").Value = 44168
$ws.Range("K13").Value = 'Castle Brite'
$ws.Range("L13").Value = 'Primera'
$ws.Range("M13").Value = 200
$ws.Range("N13").Value = 23500
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 23750
$ws.Range("Q13").Value = '$/caja 18 kilos'
$ws.Range("S13").Value = 1319
$ws.Range("T13").Value = 18
# Row 16
$ws.Range("D16").Value = 44167
$ws.Range("K16").Value = 'Castle Brite'
$ws.Range("M16").Value = 400
$ws.Range("N16").Value = 20000
$ws.Range("O16").Value = 21000
$ws.Range("P16").Value = 20500
$ws.Range("Q16").Value = '$/caja 15 kilos'
$ws.Range("S16").Value = 1367
$ws.Range("T16").Value = 15
# Row 17
$ws.Range("D17").Value = 44167
$ws.Range("M17").Value = 360
$ws.Range("N17").Value = 17000
$ws.Range("P17").Value = 17500
$ws.Range("Q17").Value = '$/caja 15 kilos'
$ws.Range("S17").Value = 1167
$ws.Range("T17").Value = 15
